$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.087.87'
$ws.Range("E2").Value = '  +4.21%  '
$ws.Range("D3").Value = '2.433.10'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.30'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +8.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.517'
$ws.Range("D7").ClearFormats()
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +9.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.83'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0805'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.122'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.97'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").Value = '2.808.62'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '2.464.56'
$ws.Range("E16").Value = '  +4.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.836'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.94%  '
$ws.Range("D18").Value = '44.973.84'
$ws.Range("E18").Value = '  +4.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.39'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.36'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '0.0₃0921'
$ws.Range("E21").Value = '  +3.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.01'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.46'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.18%  '
$ws.Range("E24").Value = '  +3.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.44'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("E28").Value = '  +10.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.60'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.28%  '
$ws.Range("E32").Value = '  +16.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.95'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +11.80%  '
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E37").Value = '  +3.79%  '
$ws.Range("E38").Value = '  +3.80%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '127.12'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.19'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").Value = '1.947.11'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.98'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +7.70%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.12'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.24'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.79'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +17.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.89'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.03'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.47%  '
